$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.658.84'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.278.62'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.92'
$ws.Range("E5").Value = '  -4.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.97'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.37'
$ws.Range("E10").Value = '  -3.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0894'
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.963'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.05'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '2.622.28'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '2.280.09'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '42.251.48'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.27'
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.37'
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.35'
$ws.Range("E23").Value = '  -7.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.13'
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.33'
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.89'
$ws.Range("E29").Value = '  +13.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.93'
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.68'
$ws.Range("E31").Value = '  -5.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.10'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0844'
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.110'
$ws.Range("E36").Value = '  -3.70%  '
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0344'
$ws.Range("E38").Value = '  -2.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.61'
$ws.Range("E39").Value = '  -2.38%  '
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.55'
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.72'
$ws.Range("E42").Value = '  +8.51%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.29'
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.83'
$ws.Range("E46").Value = '  -4.40%  '
$ws.Range("D47").Value = '1.697.99'
$ws.Range("E47").Value = '  +6.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '109.36'
$ws.Range("E48").Value = '  -3.34%  '
$ws.Range("E49").Value = '  -4.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.58'
$ws.Range("E50").Value = '  -3.62%  '
$ws.Range("E51").Value = '  -2.14%  '
